# Fruta / hortaliza, semanal
# Insert a new weekly record as row 116 (pushing the existing rows 116-122 down to 117-123),
# matching the new row's data to the style/shape of its former neighbour (old row 116).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 116; Excel shifts rows 116:122 down to 117:123
# and the inserted row inherits formatting (incl. the date style) from the row above it.
$ws.Rows("116").Insert()

# Populate the newly inserted row 116 with the new weekly record.
$ws.Range("A116").Value = 7
$ws.Range("B116").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C116").Value = "Ñuble"
$ws.Range("D116").Value = 44931
$ws.Range("E116").Value = 16
$ws.Range("F116").Value = 100112021
$ws.Range("G116").Value = "Ají"
$ws.Range("H116").Value = "Americana (o)"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 100
$ws.Range("K116").Value = 13000
$ws.Range("L116").Value = 14000
$ws.Range("M116").Value = 13500
$ws.Range("N116").Value = "$/caja 15 kilos"
$ws.Range("O116").Value = "Región del Maule"
$ws.Range("P116").Value = 900
$ws.Range("Q116").Value = 15
$ws.Range("R116").Value = "Hortaliza"

# Make sure the new date cell keeps the same date number format as the rest of column D.
$ws.Range("D116").NumberFormat = $ws.Range("D117").NumberFormat
